$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new header columns: J1 "nextkin", K1 "kinphone" ---
$ws.Range("J1").Value = "nextkin"
$ws.Range("K1").Value = "kinphone"

# Style the new header cells like the other bold header cells (bold font + new fill color)
$headerRange = $ws.Range("J1:K1")
$headerRange.Font.Bold = $true
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Interior.ThemeColor = 7
$headerRange.Borders.LineStyle = 1

# --- Update selection (active cell) to match diff ---
$ws.Range("G17").Select()
